$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top to hold the header labels; existing data shifts down.
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "year"
$ws.Range("B1").Value = "value"

$ws.Range("B1").Select()
